# Filter for account bills skipped in budget bills
# - "Prise" header/column is renamed to "Amount" (the old mis-spelled
#   "status"-adjacent "Prise" header text is no longer needed / replaced).
# - Row 4's transaction date moves forward a day (2022-01-04 -> 2022-01-05,
#   serials 44565 -> 44566) to line up with the new budget-bills fixture.
# - Selection cursor left on F1 (the renamed "Amount" header) instead of D2.
# - Page is set up for A4 portrait printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Rename the "Prise" header (column F) to "Amount".
$ws.Range("F1").Value = "Amount"

# Bump the date in row 4 by one day.
$ws.Range("D4").Value = 44566

# Page setup: A4, portrait (matches the printed layout used for the tests).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the cursor on the renamed "Amount" header cell.
$ws.Range("F1").Select()
